$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the same cell formatting as the last existing data row (200) to the
# five new rows, cell-by-cell so that columns left blank in row 200 (T, U, W)
# stay blank in the new rows too.
$styleCols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","V","X","Y","Z","AA","AB","AC","AD","AE","AF")
foreach ($newRow in 201..205) {
    foreach ($c in $styleCols) {
        $ws.Range($c + "200").Copy()
        $ws.Range($c + $newRow).PasteSpecial(-4122)
    }
}
$excel.CutCopyMode = 0

# Populate the new rows with the book data (columns A..AF, left to right).
# Row 201
$ws.Range("A201").Value = 200.0
$ws.Range("B201").Value = "Book"
$ws.Range("C201").Value = "Will"
$ws.Range("D201").Value = "Yes"
$ws.Range("E201").Value = "Will Smith"
$ws.Range("F201").Value = "Biography Of Will Smith"
$ws.Range("G201").Value = "Focus Area"
$ws.Range("H201").Value = 418.0
$ws.Range("I201").Value = "NA"
$ws.Range("J201").Value = "Rs. 239"
$ws.Range("K201").Value = 1.0
$ws.Range("L201").Value = "1-1-D"
$ws.Range("M201").Value = "Paperback"
$ws.Range("N201").Value = 2021.0
$ws.Range("O201").Value = "English"
$ws.Range("P201").Value = "Self Help"
$ws.Range("Q201").Value = "Biography"
$ws.Range("R201").Value = "No"
$ws.Range("S201").Value = "NF"
$ws.Range("V201").Value = 0.0
$ws.Range("X201").Value = "Male"
$ws.Range("Y201").Value = "978-1-529-12416-3"
$ws.Range("Z201").Value = 1.0
$ws.Range("AA201").Value = "First Floor"
$ws.Range("AB201").Value = "Practical"
$ws.Range("AC201").Value = "Yes"
$ws.Range("AD201").Value = "Yes"
$ws.Range("AE201").Value = 9.5
$ws.Range("AF201").Value = "New"

# Row 202
$ws.Range("A202").Value = 201.0
$ws.Range("B202").Value = "Book"
$ws.Range("C202").Value = "Hooked"
$ws.Range("D202").Value = "Yes"
$ws.Range("E202").Value = "Nir Eyal"
$ws.Range("F202").Value = "How to make habit forming products"
$ws.Range("G202").Value = "Focus Area"
$ws.Range("H202").Value = 246.0
$ws.Range("I202").Value = "Matt Mullenweg , Mave McLure"
$ws.Range("J202").Value = "Rs. 169"
$ws.Range("K202").Value = 2.0
$ws.Range("L202").Value = "1-1-D"
$ws.Range("M202").Value = "HardCover"
$ws.Range("N202").Value = 2019.0
$ws.Range("O202").Value = "English"
$ws.Range("P202").Value = "Business"
$ws.Range("Q202").Value = "Entreprunership"
$ws.Range("R202").Value = "No"
$ws.Range("S202").Value = "NF"
$ws.Range("V202").Value = 0.0
$ws.Range("X202").Value = "Male"
$ws.Range("Y202").Value = "978-0-241-184483-7"
$ws.Range("Z202").Value = 1.0
$ws.Range("AA202").Value = "First Floor"
$ws.Range("AB202").Value = "Practical"
$ws.Range("AC202").Value = "Yes"
$ws.Range("AD202").Value = "Yes"
$ws.Range("AE202").Value = 9.5
$ws.Range("AF202").Value = "New"

# Row 203
$ws.Range("A203").Value = 202.0
$ws.Range("B203").Value = "Book"
$ws.Range("C203").Value = "The Obstacle Is The Way"
$ws.Range("D203").Value = "Yes"
$ws.Range("E203").Value = "Ryan Holiday"
$ws.Range("F203").Value = "Turning Adversity Into Advantage"
$ws.Range("G203").Value = "Focus Area"
$ws.Range("H203").Value = 201.0
$ws.Range("I203").Value = "Robert Greene"
$ws.Range("J203").Value = "Rs. 139"
$ws.Range("K203").Value = 1.0
$ws.Range("L203").Value = "1-1-D"
$ws.Range("M203").Value = "Paperback"
$ws.Range("N203").Value = 2014.0
$ws.Range("O203").Value = "English"
$ws.Range("P203").Value = "Self Help"
$ws.Range("Q203").Value = "Power"
$ws.Range("R203").Value = "No"
$ws.Range("S203").Value = "NF"
$ws.Range("V203").Value = 0.0
$ws.Range("X203").Value = "Male"
$ws.Range("Y203").Value = "978-1-7812-5148-5"
$ws.Range("Z203").Value = 1.0
$ws.Range("AA203").Value = "First Floor"
$ws.Range("AB203").Value = "Practical"
$ws.Range("AC203").Value = "Yes"
$ws.Range("AD203").Value = "Yes"
$ws.Range("AE203").Value = 9.0
$ws.Range("AF203").Value = "New"

# Row 204
$ws.Range("A204").Value = 203.0
$ws.Range("B204").Value = "Book"
$ws.Range("C204").Value = "The Courage To Be Disliked"
$ws.Range("D204").Value = "Yes"
$ws.Range("E204").Value = "Ichiro Kishimi"
$ws.Range("F204").Value = "How to free yourself and achieve real happiness"
$ws.Range("G204").Value = "Focus Area"
$ws.Range("H204").Value = 272.0
$ws.Range("I204").Value = "NA"
$ws.Range("J204").Value = "Rs. 159"
$ws.Range("K204").Value = 1.0
$ws.Range("L204").Value = "1-1-D"
$ws.Range("M204").Value = "Paperback"
$ws.Range("N204").Value = 2018.0
$ws.Range("O204").Value = "English"
$ws.Range("P204").Value = "Self Help"
$ws.Range("Q204").Value = "Psychology"
$ws.Range("R204").Value = "No"
$ws.Range("S204").Value = "NF"
$ws.Range("V204").Value = 0.0
$ws.Range("X204").Value = "Male"
$ws.Range("Y204").Value = "978-1-76063-072-0"
$ws.Range("Z204").Value = 1.0
$ws.Range("AA204").Value = "First Floor"
$ws.Range("AB204").Value = "Practical"
$ws.Range("AC204").Value = "Yes"
$ws.Range("AD204").Value = "Yes"
$ws.Range("AE204").Value = 9.0
$ws.Range("AF204").Value = "New"

# Row 205
$ws.Range("A205").Value = 204.0
$ws.Range("B205").Value = "Book"
$ws.Range("C205").Value = "the 5 Second Rule"
$ws.Range("D205").Value = "Yes"
$ws.Range("E205").Value = "Mel Robbins"
$ws.Range("F205").Value = "Transform Your Life With Everyday Courage"
$ws.Range("G205").Value = "Focus Area"
$ws.Range("H205").Value = 238.0
$ws.Range("I205").Value = "NA"
$ws.Range("J205").Value = "Rs. 210"
$ws.Range("K205").Value = 1.0
$ws.Range("L205").Value = "1-1-D"
$ws.Range("M205").Value = "HardCover"
$ws.Range("N205").Value = 2017.0
$ws.Range("O205").Value = "English"
$ws.Range("P205").Value = "Self Help"
$ws.Range("Q205").Value = "Psychology"
$ws.Range("R205").Value = "No"
$ws.Range("S205").Value = "NF"
$ws.Range("V205").Value = 0.0
$ws.Range("X205").Value = "Female"
$ws.Range("Y205").Value = "978-1-68261-238-5"
$ws.Range("Z205").Value = 1.0
$ws.Range("AA205").Value = "First Floor"
$ws.Range("AB205").Value = "Practical"
$ws.Range("AC205").Value = "Yes"
$ws.Range("AD205").Value = "Yes"
$ws.Range("AE205").Value = 9.8
$ws.Range("AF205").Value = "New"

